# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The sheet lists worker payment-arrears ("estado de cuenta"). This edit
# replaces the worker/period rows with a new data set (only 2 workers /
# 2 periods now, instead of 4 workers / 8 periods), updates the summary
# totals accordingly, and removes the now-unused rows so the "firma"
# block moves back up right under the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: first worker record -> CARLOS JAIR GONZALEZ ROMERO ----------
$ws.Range("C16").Value = "20111995"
$ws.Range("D16").Value = "CARLOS JAIR GONZALEZ ROMERO"
$ws.Range("E16").Value = "1712"
$ws.Range("F16").Value = 14486
$ws.Range("G16").Value = 781242

# --- Row 17: second worker record -> JERGES DAVID CASTRO ROMERO ----------
$ws.Range("C17").Value = "73209651"
$ws.Range("D17").Value = "JERGES DAVID CASTRO ROMERO"
$ws.Range("E17").Value = "2107"
$ws.Range("F17").Value = 203220
$ws.Range("G17").Value = 9175200

# Row 17 becomes the new final data row, so it should pick up the bolder
# "subtotal" formatting that used to belong to the old last row (29),
# before that row gets removed below.
$ws.Range("B29:J29").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Summary header values -------------------------------------------------
$ws.Range("E11").Value = 217706   # VALOR MORA total
$ws.Range("C13").Value = 2        # Cant. Trabajadores
$ws.Range("F13").Value = 2        # Cant. Periodos

# --- Remove the now-obsolete worker/period rows and the old totals row ---
# (old rows 18-29: extra periods for worker 1, the whole worker "EDHILBERTO"
# block, and the old CARLOS/totals row -- all superseded by rows 16-17
# above). Deleting these entire rows shifts the signature block (old rows
# 34-35) up to rows 22-23, matching the original gap of 4 blank rows.
$ws.Range("B18:J29").EntireRow.Delete()
